# Auto update Excel log
# Appends newly captured sensor readings (2026-01-28, ~12:04-12:17) to the
# PIR, Humidity and Temperature logs. Column A (Date) and, on the Humidity
# sheet, column E (the "NN.N%" reading) look like a date / a percentage to
# Excel's automatic type detection. Each such cell is briefly forced to
# Text format before the literal value is written, then its style is put
# back to Normal -- this keeps the value as a plain string (matching every
# other row in the log) instead of being reinterpreted as a date serial
# number or a percentage value.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 177-189 ---
$wsPir = $wb.Worksheets.Item("PIR")
$pirRows = @(
  @(177, '2026-01-28', '12:04:55', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(178, '2026-01-28', '12:16:32', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(179, '2026-01-28', '12:16:36', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(180, '2026-01-28', '12:16:41', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(181, '2026-01-28', '12:16:46', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(182, '2026-01-28', '12:16:52', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(183, '2026-01-28', '12:16:56', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(184, '2026-01-28', '12:17:01', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(185, '2026-01-28', '12:17:06', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(186, '2026-01-28', '12:17:12', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(187, '2026-01-28', '12:17:16', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(188, '2026-01-28', '12:17:21', '12:00', 'Bathroom', 'No Motion', 'Inactive'),
  @(189, '2026-01-28', '12:17:26', '12:00', 'Bathroom', 'No Motion', 'Inactive')
)
foreach ($r in $pirRows) {
    $rowNum = $r[0]

    $dateCell = $wsPir.Range("A$rowNum")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r[1]
    $dateCell.Style = "Normal"

    $wsPir.Range("B$rowNum").Value = $r[2]
    $wsPir.Range("C$rowNum").Value = $r[3]
    $wsPir.Range("D$rowNum").Value = $r[4]
    $wsPir.Range("E$rowNum").Value = $r[5]
    $wsPir.Range("F$rowNum").Value = $r[6]
}

# --- Humidity sheet: append rows 166-178 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
  @(166, '2026-01-28', '12:04:54', '12:00', 'Bathroom', '88.1%', 'Active'),
  @(167, '2026-01-28', '12:16:31', '12:00', 'Bathroom', '87.5%', 'Active'),
  @(168, '2026-01-28', '12:16:35', '12:00', 'Bathroom', '86.6%', 'Active'),
  @(169, '2026-01-28', '12:16:39', '12:00', 'Bathroom', '87.6%', 'Active'),
  @(170, '2026-01-28', '12:16:43', '12:00', 'Bathroom', '87.6%', 'Active'),
  @(171, '2026-01-28', '12:16:51', '12:00', 'Bathroom', '87.6%', 'Active'),
  @(172, '2026-01-28', '12:16:55', '12:00', 'Bathroom', '86.7%', 'Active'),
  @(173, '2026-01-28', '12:16:59', '12:00', 'Bathroom', '87.6%', 'Active'),
  @(174, '2026-01-28', '12:17:03', '12:00', 'Bathroom', '87.6%', 'Active'),
  @(175, '2026-01-28', '12:17:07', '12:00', 'Bathroom', '86.7%', 'Active'),
  @(176, '2026-01-28', '12:17:11', '12:00', 'Bathroom', '87.7%', 'Active'),
  @(177, '2026-01-28', '12:17:15', '12:00', 'Bathroom', '87.4%', 'Active'),
  @(178, '2026-01-28', '12:17:23', '12:00', 'Bathroom', '88.0%', 'Active')
)
foreach ($r in $humidityRows) {
    $rowNum = $r[0]

    $dateCell = $wsHumidity.Range("A$rowNum")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r[1]
    $dateCell.Style = "Normal"

    $wsHumidity.Range("B$rowNum").Value = $r[2]
    $wsHumidity.Range("C$rowNum").Value = $r[3]
    $wsHumidity.Range("D$rowNum").Value = $r[4]

    $valueCell = $wsHumidity.Range("E$rowNum")
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $r[5]
    $valueCell.Style = "Normal"

    $wsHumidity.Range("F$rowNum").Value = $r[6]
}

# --- Temperature sheet: append rows 166-178 ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
  @(166, '2026-01-28', '12:04:55', '12:00', 'Bathroom', '23.0C', 'Active'),
  @(167, '2026-01-28', '12:16:31', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(168, '2026-01-28', '12:16:35', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(169, '2026-01-28', '12:16:39', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(170, '2026-01-28', '12:16:43', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(171, '2026-01-28', '12:16:52', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(172, '2026-01-28', '12:16:56', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(173, '2026-01-28', '12:17:00', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(174, '2026-01-28', '12:17:04', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(175, '2026-01-28', '12:17:08', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(176, '2026-01-28', '12:17:12', '12:00', 'Bathroom', '22.9C', 'Active'),
  @(177, '2026-01-28', '12:17:16', '12:00', 'Bathroom', '23.0C', 'Active'),
  @(178, '2026-01-28', '12:17:24', '12:00', 'Bathroom', '23.0C', 'Active')
)
foreach ($r in $temperatureRows) {
    $rowNum = $r[0]

    $dateCell = $wsTemperature.Range("A$rowNum")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r[1]
    $dateCell.Style = "Normal"

    $wsTemperature.Range("B$rowNum").Value = $r[2]
    $wsTemperature.Range("C$rowNum").Value = $r[3]
    $wsTemperature.Range("D$rowNum").Value = $r[4]
    $wsTemperature.Range("E$rowNum").Value = $r[5]
    $wsTemperature.Range("F$rowNum").Value = $r[6]
}
